# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (TB), C (d2S), D (K), E (IP), G (sum)
# F (Win) is left unchanged.

$data = @{
    2 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.1496068669990043; E = 0.5333859586016987;  G = 5.582307763322248 }
    3 = @{ B = 0.6545652718822623; C = 1.626987699542094;  D = 0.7210945179870265; E = 13.86384647080068;   G = 16.86649396021207 }
    4 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
    5 = @{ B = 3.272327238179451;  C = 0.3048912486333797; D = 0.7210945179870265; E = 0.5333859586016987;  G = 4.831698963401555 }
    6 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 4.327115817150455 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
